$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Widen column F (6th column) from 18 to 26
# (25.14 is the ColumnWidth input value that serializes to exactly width="26" in the OOXML)
$ws.Columns.Item(6).ColumnWidth = 25.14

# Update VENTA (D3), POR CUMPLIR (E3) and CUMPLIMIENTO (F3) for row 3 (PORCELANATO)
$ws.Cells.Item(3, 4).Value = 24.31
$ws.Cells.Item(3, 5).Value = 17475.69
$ws.Cells.Item(3, 6).Value = 0.001389142857142857

# Update the TOTAL row (row 4) to reflect the new sums
$ws.Cells.Item(4, 4).Value = 772.3499999999999
$ws.Cells.Item(4, 5).Value = 16727.65
$ws.Cells.Item(4, 6).Value = 0.04413428571428571
